$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48

$ws.Cells.Item($row, 1).Value = 45952
$ws.Cells.Item($row, 2).Value = "21,7178"
$ws.Cells.Item($row, 3).Value = "15,6198"
$ws.Cells.Item($row, 4).Value = "15,4273"
$ws.Cells.Item($row, 5).Value = "15,4273"

$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
